# Edit script: "Unveiling the Mysteries of Dark Matter" -> "Biology: The Tapestry of Life"
# Applies the font-name fix (TimesNewToman -> Times New Roman) and swaps the
# dark-matter themed content for biology themed content, including two new
# sentences appended within the body paragraph and one new sentence appended
# within the summary paragraph, plus a trailing empty paragraph.

$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Fix the (misspelled) font name across the whole document body.
# ---------------------------------------------------------------------------
$full = $d.Range(0, $d.Content.End)
$full.Font.Name = "Times New Roman"

# ---------------------------------------------------------------------------
# 2) Straightforward text substitutions (one run's text for another).
# ---------------------------------------------------------------------------
function Replace-Text($oldText, $newText) {
    $ok = $d.Content.Find.Execute($oldText, $true, $false, $false, $false, $false, `
                                   $true, 1, $false, $newText, 2)
    if (-not $ok) {
        throw "Find/Replace failed for: $oldText"
    }
}

Replace-Text "Unveiling the Mysteries of Dark Matter" "Biology: The Tapestry of Life"
Replace-Text "Alana Hendricks" "Sarah Johnson"
Replace-Text "alanahendricks@astronomysociety" "sarahj@educatormail"

Replace-Text "For decades, astronomers and physicists have been engrossed in a captivating cosmic enigma: the existence and nature of dark matter" `
             "Biology, the study of life, is an awe-inspiring exploration into the diverse tapestry of organisms that inhabit our planet"

Replace-Text "This enigmatic substance, believed to comprise over 26% of the universe, exerts a gravitational influence far exceeding its apparent mass, shaping galaxies, and influencing the universe's expansion rate" `
             "It captivates us with its complexities and unravels the intricate web that connects all living things"

Replace-Text "Yet, despite its profound impact, dark matter remains shrouded in mystery, fueling scientific curiosity and propelling groundbreaking research" `
             "Embark on a journey to discover the vibrant ecosystems that thrive in harmony, the astonishing adaptations that organisms possess, and the profound impact we have on the delicate balance of life"

Replace-Text "The evidence for dark matter's existence is compelling" `
             "As we delve into the microscopic realm of cells, we uncover astonishing marvels of organization and functionality"

Replace-Text "Through meticulous observations of galaxies and galaxy clusters, astronomers have discovered that the gravitational force necessary to hold these celestial structures together far surpasses the gravitational pull exerted by the visible matter they contain" `
             "Each cell is an intricate microcosm within itself, carrying out essential processes that sustain life"

Replace-Text "This discrepancy suggests the presence of an invisible mass, an unseen entity governing the universe's dynamics" `
             "The diversity of life becomes evident as we encounter the remarkable variations in form and behavior among organisms"

Replace-Text "Furthermore, observations of the cosmic microwave background radiation, the leftover glow from the Big Bang, provide further clues about dark matter's existence" `
             "Biology enables us to comprehend how organisms interact with each other and with their surroundings, forming complex ecosystems that thrive through interconnectedness"

Replace-Text "Minute temperature variations in this radiation hint at the gravitational influence of dark matter during the universe's early moments, supporting the notion that it played a pivotal role in shaping the universe's structure" `
             "We learn how delicate balances are maintained within these ecosystems and the intricate roles that each organism plays in preserving this equilibrium"

Replace-Text "The quest to understand dark matter has captivated scientists worldwide, driving cutting-edge research and groundbreaking discoveries" `
             "Biology is an intriguing subject that unveils the mysteries of life's tapestry"

Replace-Text "While its true identity remains elusive, the evidence for its existence is undeniable" `
             "It encompasses the study of cells, their intricate organization and functionality, the bewildering diversity of organisms, and the dynamic interactions between organisms within ecosystems"

Replace-Text "Through continued exploration and innovation, scientists are determined to unravel the mysteries of dark matter, shedding light on one of the universe's most enigmatic components" `
             "Biology nurtures an understanding of the influence we have on the environment and inspires us to act as responsible stewards of our planet"

# ---------------------------------------------------------------------------
# 3) Append the three brand-new sentences (each followed by its own "."
#    run, matching the surrounding document's one-sentence-per-run style).
# ---------------------------------------------------------------------------
function Append-Sentence($anchorText, $sentence) {
    $found = $d.Content.Find.Execute($anchorText)
    if (-not $found) {
        throw "Could not locate anchor text: $anchorText"
    }
    $r = $d.Content
    $r.Find.Execute($anchorText) | Out-Null
    $r.Collapse(0)
    $r.InsertAfter($sentence)
    $r.Collapse(0)
    $r.InsertAfter(".")
}

Append-Sentence "The diversity of life becomes evident as we encounter the remarkable variations in form and behavior among organisms." `
    " From the grandeur of the majestic whales that roam our oceans to the minuscule yet tenacious microorganisms, we marvel at the myriad life forms that grace our planet"

Append-Sentence "We learn how delicate balances are maintained within these ecosystems and the intricate roles that each organism plays in preserving this equilibrium." `
    " The interdependence of organisms within these interconnected systems highlights the profound responsibility we bear as stewards of our natural world"

Append-Sentence "Biology nurtures an understanding of the influence we have on the environment and inspires us to act as responsible stewards of our planet." `
    " As we continue to explore the wonders of life, we unravel the secrets of our existence and uncover the boundless possibilities that the study of biology holds"

# ---------------------------------------------------------------------------
# 4) Add the trailing empty paragraph that now follows the Summary section.
# ---------------------------------------------------------------------------
$d.Content.InsertParagraphAfter()
